# Ajout encryptedPassword + ajout autres eleves CSV
# Adds four new students (Nom, Prenom, Email) below the existing list and
# turns their Email cells into mailto: hyperlinks, matching the original
# author's manual data entry + "Insert Hyperlink" workflow in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New student rows -------------------------------------------------

$ws.Range("A6").Value = "Lorthioir"
$ws.Range("B6").Value = "Jérémy"
$ws.Range("C6").Value = "Lorthioir.Jeremy@etu.unilim.fr"

$ws.Range("A7").Value = "Sparrow"
$ws.Range("B7").Value = "Jack"
$ws.Range("C7").Value = "Sparrow.Jack@pirate.com"

$ws.Range("A8").Value = "Potter"
$ws.Range("B8").Value = "Harry"
$ws.Range("C8").Value = "HarryPotter@poudlard.com"

$ws.Range("A9").Value = "Picsou"
$ws.Range("B9").Value = "Balthazar"
$ws.Range("C9").Value = "Balthazar.Picsou@riche.com"

# --- Hyperlink-ify the new Email cells (mailto:) -----------------------

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:Lorthioir.Jeremy@etu.unilim.fr")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:Sparrow.Jack@pirate.com")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:HarryPotter@poudlard.com")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:Balthazar.Picsou@riche.com")

# --- Leave the selection where the author's cursor ended up ------------

$null = $ws.Range("A10").Select()
